# Add four new "Costsheet" rows (82-85) to Sheet1, mirroring the layout
# used by the existing RuleSet sections (Task/Organization/Person/etc.).
#
# New shared strings introduced by this edit (created implicitly, in this
# order, as the cell values below are written):
#   - Costsheet – Only participants can add files
#   - COSTSHEET
#   - grant add file to assignee, owner, co-owner, supervisor, owning group, approver, collaborator
#   - Costsheet – default read access
#   - Costsheet – Only participants can save
#   - grant save to assignee, owner, co-owner, supervisor, owning group, approver, collaborator
#   - grant uploadOrReplaceFile to assignee, owner, co-owner, supervisor, owning group, approver, collaborator
#   - Costsheet –  Only participants can upload or replace files

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The authored workbook switched the calculation mode to manual.
$xlCalculationManual = -4135
$excel.Calculation = $xlCalculationManual

$xlPasteFormats = -4122

# --- Row 83: "Costsheet – Only participants can add files" -----------
# Same visual layout as row 48 ("... – Only participants can add files" style).
$ws.Range("A48:G48").Copy() | Out-Null
$ws.Range("A83:G83").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B83").Value = "Costsheet – Only participants can add files"
$ws.Range("C83").Value = "COSTSHEET"
$ws.Range("G83").Value = "grant add file to assignee, owner, co-owner, supervisor, owning group, approver, collaborator"

# --- Row 82: "Costsheet – default read access" -----------------------
# Same visual layout as row 60 ("Folder - default list folder").
$ws.Range("A60:G60").Copy() | Out-Null
$ws.Range("A82:G82").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B82").Value = "Costsheet – default read access"
$ws.Range("C82").Value = "COSTSHEET"
$ws.Range("G82").Value = "grant read to *"

# --- Row 84: "Costsheet – Only participants can save" ----------------
$ws.Range("A48:G48").Copy() | Out-Null
$ws.Range("A84:G84").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B84").Value = "Costsheet – Only participants can save"
$ws.Range("C84").Value = "COSTSHEET"
$ws.Range("G84").Value = "grant save to assignee, owner, co-owner, supervisor, owning group, approver, collaborator"

# --- Row 85: "Costsheet –  Only participants can upload or replace files"
# Same visual layout as row 26 (column D uses the "164/wrap" style).
$ws.Range("A26:G26").Copy() | Out-Null
$ws.Range("A85:G85").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("G85").Value = "grant uploadOrReplaceFile to assignee, owner, co-owner, supervisor, owning group, approver, collaborator"
$ws.Range("B85").Value = "Costsheet –  Only participants can upload or replace files"
$ws.Range("C85").Value = "COSTSHEET"

# Restore the selection to reflect where the edit was made.
$ws.Range("C85").Select() | Out-Null

Write-Host "Added Costsheet rows 82-85 to Sheet1"
